$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 19 (Africana... / Fine Arts block),
# shifting the existing rows 19-114 down to 21-116.
$ws.Rows("19:20").Insert()

# Fill in the new rows with the Fine Arts Reserve mappings.
$ws.Range("A19").Value = "Fine Arts Library Permanent Reserve"
$ws.Range("E19").Value = "Fine Arts Library > Reserve"

$ws.Range("A20").Value = "Fine Arts Course Reserve (Ask at Circulation)"
$ws.Range("E20").Value = "Fine Arts Library > Reserve"

# Match the saved selection/view state from the edit.
$ws.Range("D20").Select()
